# DataEngine.xlsx maintenance edit:
#  - Add a new "TestCases" sheet in front of the existing "TestSteps" sheet.
#  - Populate "TestCases" with the two login test cases and a RunMode column.
#  - Rename the old generic "Successfull_login_01" TestCase_ID used throughout
#    "TestSteps" to the new "login_01" id that matches "Login_01" in TestCases.

$wb = $excel.ActiveWorkbook
$testSteps = $wb.Worksheets.Item("TestSteps")

# Update the TestCase_ID column on TestSteps: every data row (2-10) used to
# reference "Successfull_login_01"; it now points at "login_01".
$testSteps.Range("A2").Value = "login_01"
$testSteps.Range("A3").Value = "login_01"
$testSteps.Range("A4").Value = "login_01"
$testSteps.Range("A5").Value = "login_01"
$testSteps.Range("A6").Value = "login_01"
$testSteps.Range("A7").Value = "login_01"
$testSteps.Range("A8").Value = "login_01"
$testSteps.Range("A9").Value = "login_01"
$testSteps.Range("A10").Value = "login_01"

# Clear the stale selection/active cell left over on TestSteps now that
# TestCases becomes the front/selected sheet.
$testSteps.Activate()
$testSteps.Range("A1").Select()

# Insert the new "TestCases" sheet before "TestSteps" so it becomes the
# first (and now active/selected) tab in the workbook.
$testCases = $wb.Worksheets.Add($testSteps)
$testCases.Name = "TestCases"

$testCases.Range("A1").Value = "TestCase_ID"
$testCases.Range("B1").Value = "Description"
$testCases.Range("C1").Value = "RunMode"

$testCases.Range("A2").Value = "Login_01"
$testCases.Range("B2").Value = "Successful Registration on entering same password in both the fields"
$testCases.Range("C2").Value = "Yes"

$testCases.Range("A3").Value = "Login_02"
$testCases.Range("B3").Value = "Error message on entering different values in both the fields"
$testCases.Range("C3").Value = "Yes"

# Re-use the existing header / body cell formatting from TestSteps (bold
# white-on-black header row, thin-bordered body rows).
$testSteps.Range("A1:C1").Copy()
$testCases.Range("A1:C1").PasteSpecial(-4122)

$testSteps.Range("A2:C2").Copy()
$testCases.Range("A2:C3").PasteSpecial(-4122)

$testCases.Columns.Item(1).ColumnWidth = 19.65
$testCases.Columns.Item(2).ColumnWidth = 62.65

$testCases.Range("A2:C3").Select()
